# Update build timestamp references from "17.29.55 EST" to "18.05.36 EST"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 13; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value()
    $cell.Value = $val.Replace($oldStamp, $newStamp)
}
